$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "asistencia"

# --- New weekly attendance columns AE (08/09) and AF (15/09) ---
$dateFormat = $ws.Range("E3").NumberFormat

$ws.Range("AE3").Value = 45908
$ws.Range("AE3").NumberFormat = $dateFormat
$ws.Range("AF3").Value = 45915
$ws.Range("AF3").NumberFormat = $dateFormat

# Attendance values per student row (row 4..29), AE = 08/09, AF = 15/09
$ae = @("P","P","A","P","P","A","P","P","P","P","P","P","P","P","P","P","P","P","P","P","P","P","A","P","P","P")
$af = @("P","P","P","P","P","P","P","P","P","P","P","P","A","P","P","P","P","P","P","A","P","P","P","O","P","P")

for ($i = 0; $i -lt 26; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 31).Value = $ae[$i]   # column AE = 31
    $ws.Cells.Item($row, 32).Value = $af[$i]   # column AF = 32
}

# Totals row
$ws.Range("AE30").Formula = '=COUNTIF(AE4:AE29,"P")'
$ws.Range("AF30").Formula = '=COUNTIF(AF4:AF29,"P")'

# --- View state: asistencia becomes the active tab/sheet (was "positivos") ---
$ws.Activate()
$excel.Goto($ws.Range("A15"), $true)
$ws.Range("AE30").Select()
